# This workbook's rows 16-21 were reordered (a permutation of whole rows).
# Mapping (new row <- old row):
#   16 <- 17
#   17 <- 18
#   18 <- 20
#   19 <- 21
#   20 <- 16
#   21 <- 19
#
# We stage the original rows in a scratch area far below the data (rows 100-105)
# and then move them into their final destination. This two-phase approach avoids
# overwriting source rows before they have been read, since the row mapping is a
# permutation (cycles), not a simple shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Phase 1: stash current rows 16-21 into scratch rows 100-105 ----
$ws.Range("A16:AY16").Cut($ws.Range("A100:AY100"))
$ws.Range("A17:AY17").Cut($ws.Range("A101:AY101"))
$ws.Range("A18:AY18").Cut($ws.Range("A102:AY102"))
$ws.Range("A19:AY19").Cut($ws.Range("A103:AY103"))
$ws.Range("A20:AY20").Cut($ws.Range("A104:AY104"))
$ws.Range("A21:AY21").Cut($ws.Range("A105:AY105"))

# ---- Phase 2: move scratch rows into their final destination ----
# new row 16 <- old row 17 (scratch 101)
$ws.Range("A101:AY101").Cut($ws.Range("A16:AY16"))
# new row 17 <- old row 18 (scratch 102)
$ws.Range("A102:AY102").Cut($ws.Range("A17:AY17"))
# new row 18 <- old row 20 (scratch 104)
$ws.Range("A104:AY104").Cut($ws.Range("A18:AY18"))
# new row 19 <- old row 21 (scratch 105)
$ws.Range("A105:AY105").Cut($ws.Range("A19:AY19"))
# new row 20 <- old row 16 (scratch 100)
$ws.Range("A100:AY100").Cut($ws.Range("A20:AY20"))
# new row 21 <- old row 19 (scratch 103)
$ws.Range("A103:AY103").Cut($ws.Range("A21:AY21"))

# ---- Phase 3: clean up phantom empty cells that Cut() creates for columns ----
# that had no cell at all in the source row (so the output keeps the same
# sparse cell layout as before, instead of materializing empty cells for
# every column in the A:AY span).
$ws.Range("J16:O16").ClearContents()
$ws.Range("X16").ClearContents()
$ws.Range("AF16").ClearContents()
$ws.Range("AH16:AS16").ClearContents()
$ws.Range("AU16:AV16").ClearContents()

$ws.Range("J17:O17").ClearContents()
$ws.Range("X17").ClearContents()
$ws.Range("AC17").ClearContents()
$ws.Range("AF17").ClearContents()
$ws.Range("AH17:AS17").ClearContents()
$ws.Range("AU17:AV17").ClearContents()

$ws.Range("J18").ClearContents()
$ws.Range("O18").ClearContents()
$ws.Range("X18").ClearContents()
$ws.Range("AC18").ClearContents()
$ws.Range("AF18").ClearContents()
$ws.Range("AH18:AS18").ClearContents()
$ws.Range("AU18:AV18").ClearContents()

$ws.Range("J19").ClearContents()
$ws.Range("O19").ClearContents()
$ws.Range("X19").ClearContents()
$ws.Range("AF19").ClearContents()
$ws.Range("AH19:AS19").ClearContents()
$ws.Range("AU19:AV19").ClearContents()

$ws.Range("J20:O20").ClearContents()
$ws.Range("X20").ClearContents()
$ws.Range("AF20").ClearContents()
$ws.Range("AH20:AS20").ClearContents()
$ws.Range("AU20:AV20").ClearContents()

$ws.Range("J21").ClearContents()
$ws.Range("O21").ClearContents()
$ws.Range("X21").ClearContents()
$ws.Range("AF21").ClearContents()
$ws.Range("AH21:AS21").ClearContents()
$ws.Range("AU21:AV21").ClearContents()
